$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G15").Formula = "=G17/4"
$ws.Range("G15").NumberFormat = "0.00000000"
Write-Host "Formula G15:" $ws.Range("G15").Formula
Write-Host "Value2 G15:" $ws.Range("G15").Value2
Write-Host "NumberFormat G15:" $ws.Range("G15").NumberFormat
